$wb = $excel.ActiveWorkbook
$ws4 = $wb.Worksheets.Item("Sheet3")
$ws4.Delete()
foreach ($s in $wb.Worksheets) {
    Write-Host $s.Name
}
